$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$dates = @("10-09-2021", "13-09-2021", "14-09-2021", "15-09-2021", "16-09-2021", "20-09-2021")

$startRow = 176
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $cellA = $ws.Cells.Item($row, 1)
    # Assign as a text formula first so Excel doesn't auto-convert the
    # dd-mm-yyyy-looking string into a date serial number, then convert
    # the formula result into a plain static value (keeps default style).
    $cellA.Formula = '="' + $dates[$i] + '"'
    $cellA.Copy()
    $cellA.PasteSpecial(-4163)
    $ws.Cells.Item($row, 2).Value = 3.25
}
